$wb = $excel.ActiveWorkbook

# ----- Sheet3 (債務) -----
$ws3 = $wb.Worksheets.Item(3)

# Row 1 headers
$ws3.Range("C1").Value = "debtor"
$ws3.Range("D1").Value = "owner"
$ws3.Range("E1").Value = "total"
$ws3.Range("F1").Value = "register_date"
$ws3.Range("G1").Value = "register_reason"
$ws3.Range("H1").Value = "property_category"
$ws3.Range("I1").Value = "category"
$ws3.Range("J1").Value = "date"
$ws3.Range("K1").Value = "legislator_name"
$ws3.Range("L1").Value = "legislator_id"
$ws3.Range("M1").Value = "source_file"
$ws3.Range("N1").Value = "index"

# Row 2
$ws3.Range("B2").Value = "催收款"
$ws3.Range("C2").Value = "顏清標"
$ws3.Range("D2").Value = "國泰世華銀行臺北市信義區松疒路"
$ws3.Range("E2").Value = 5707475
$ws3.Range("F2").Value = "85年05月30日"
$ws3.Range("G2").Value = "短期放款"
$ws3.Range("H2").Value = "debt"
$ws3.Range("I2").Value = "normal"
$ws3.Range("J2").Value = "2012-04-10"
$ws3.Range("K2").Value = "顏清標"
$ws3.Range("L2").Value = 979
$ws3.Range("M2").Value = "tmp1b4d1"
$ws3.Range("N2").Value = 87

# Row 3
$ws3.Range("B3").Value = "長期擔保放款"
$ws3.Range("C3").Value = "顏清標"
$ws3.Range("D3").Value = "華南銀行清水分行臺中市清水區中山路"
$ws3.Range("E3").Value = 8646484
$ws3.Range("F3").Value = "84年07月08日"
$ws3.Range("G3").Value = "繼承債務"
$ws3.Range("H3").Value = "debt"
$ws3.Range("I3").Value = "normal"
$ws3.Range("J3").Value = "2012-04-10"
$ws3.Range("K3").Value = "顏清標"
$ws3.Range("L3").Value = 979
$ws3.Range("M3").Value = "tmp1b4d1"
$ws3.Range("N3").Value = 89

# Row 4
$ws3.Range("B4").Value = "長期擔保放款"
$ws3.Range("C4").Value = "顔清標"
$ws3.Range("D4").Value = "華南銀行清水分行臺中市清水區中山路"
$ws3.Range("E4").Value = 9601734
$ws3.Range("F4").Value = "84年07月08日"
$ws3.Range("G4").Value = "繼承保證債務"
$ws3.Range("H4").Value = "debt"
$ws3.Range("I4").Value = "normal"
$ws3.Range("J4").Value = "2012-04-10"
$ws3.Range("K4").Value = "顏清標"
$ws3.Range("L4").Value = 979
$ws3.Range("M4").Value = "tmp1b4d1"
$ws3.Range("N4").Value = 90

# Row 5
$ws3.Range("B5").Value = "擔保放款"
$ws3.Range("C5").Value = "顔清標"
$ws3.Range("D5").Value = "元營建設股份有限公司臺中市沙鹿區北勢東路"
$ws3.Range("E5").Value = 71062315
$ws3.Range("F5").Value = "97年10月23日"
$ws3.Range("G5").Value = "依據台灣台北地方法院97年10月23H北院隆97執"
$ws3.Range("H5").Value = "debt"
$ws3.Range("I5").Value = "normal"
$ws3.Range("J5").Value = "2012-04-10"
$ws3.Range("K5").Value = "顏清標"
$ws3.Range("L5").Value = 979
$ws3.Range("M5").Value = "tmp1b4d1"
$ws3.Range("N5").Value = 91

# Row 6
$ws3.Range("B6").Value = "長期擔保放款"
$ws3.Range("C6").Value = "黃美貴"
$ws3.Range("D6").Value = "華南銀行清水分行臺中市清水區中山路"
$ws3.Range("E6").Value = 8961620
$ws3.Range("F6").Value = "84年07月08日"
$ws3.Range("G6").Value = "房屋貸款"
$ws3.Range("H6").Value = "debt"
$ws3.Range("I6").Value = "normal"
$ws3.Range("J6").Value = "2012-04-10"
$ws3.Range("K6").Value = "顏清標"
$ws3.Range("L6").Value = 979
$ws3.Range("M6").Value = "tmp1b4d1"
$ws3.Range("N6").Value = 92

# Row 7
$ws3.Range("B7").Value = "催收款"
$ws3.Range("C7").Value = "顔清標"
$ws3.Range("D7").Value = "國泰世華銀行臺北市信義區松仁路"
$ws3.Range("E7").Value = 29291128
$ws3.Range("F7").Value = "85年05月30日"
$ws3.Range("G7").Value = "短期放款"
$ws3.Range("H7").Value = "debt"
$ws3.Range("I7").Value = "normal"
$ws3.Range("J7").Value = "2012-04-10"
$ws3.Range("K7").Value = "顏清標"
$ws3.Range("L7").Value = 979
$ws3.Range("M7").Value = "tmp1b4d1"
$ws3.Range("N7").Value = 93

# Row 8
$ws3.Range("B8").Value = "催收款"
$ws3.Range("C8").Value = "顔清標"
$ws3.Range("D8").Value = "國泰世華銀行臺北市信義區松仁路"
$ws3.Range("E8").Value = 25934529
$ws3.Range("F8").Value = "85年05月30日"
$ws3.Range("G8").Value = "短期放款"
$ws3.Range("H8").Value = "debt"
$ws3.Range("I8").Value = "normal"
$ws3.Range("J8").Value = "2012-04-10"
$ws3.Range("K8").Value = "顏清標"
$ws3.Range("L8").Value = 979
$ws3.Range("M8").Value = "tmp1b4d1"
$ws3.Range("N8").Value = 94

# Row 9
$ws3.Range("B9").Value = "催收款"
$ws3.Range("C9").Value = "顏清標"
$ws3.Range("D9").Value = "國泰世華銀行臺北市信義區松仁路"
$ws3.Range("E9").Value = 24224554
$ws3.Range("F9").Value = "85年05月30日"
$ws3.Range("G9").Value = "短期放款"
$ws3.Range("H9").Value = "debt"
$ws3.Range("I9").Value = "normal"
$ws3.Range("J9").Value = "2012-04-10"
$ws3.Range("K9").Value = "顏清標"
$ws3.Range("L9").Value = 979
$ws3.Range("M9").Value = "tmp1b4d1"
$ws3.Range("N9").Value = 95

# Row 10
$ws3.Range("B10").Value = "催收款"
$ws3.Range("C10").Value = "顏清標"
$ws3.Range("D10").Value = "國泰世華銀行臺北市信義區松仁路"
$ws3.Range("E10").Value = 5742274
$ws3.Range("F10").Value = "85年11月16H"
$ws3.Range("G10").Value = "房屋貸款"
$ws3.Range("H10").Value = "debt"
$ws3.Range("I10").Value = "normal"
$ws3.Range("J10").Value = "2012-04-10"
$ws3.Range("K10").Value = "顏清標"
$ws3.Range("L10").Value = 979
$ws3.Range("M10").Value = "tmp1b4d1"
$ws3.Range("N10").Value = 96

# Copy header style onto new header cells
$ws3.Range("G1").Copy()
$ws3.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Copy data row style onto new data cells (best-effort)
$ws3.Range("G2").Copy()
$ws3.Range("H2:N2").PasteSpecial(-4122)
$ws3.Range("G3").Copy()
$ws3.Range("H3:N3").PasteSpecial(-4122)
$ws3.Range("G4").Copy()
$ws3.Range("H4:N4").PasteSpecial(-4122)
$ws3.Range("G5").Copy()
$ws3.Range("H5:N5").PasteSpecial(-4122)
$ws3.Range("G6").Copy()
$ws3.Range("H6:N6").PasteSpecial(-4122)
$ws3.Range("G7").Copy()
$ws3.Range("H7:N7").PasteSpecial(-4122)
$ws3.Range("G8").Copy()
$ws3.Range("H8:N8").PasteSpecial(-4122)
$ws3.Range("G9").Copy()
$ws3.Range("H9:N9").PasteSpecial(-4122)
$ws3.Range("G10").Copy()
$ws3.Range("H10:N10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ----- Sheet4 (事業投資) -----
$ws4 = $wb.Worksheets.Item(4)

# Row 1 headers
$ws4.Range("B1").Value = "owner"
$ws4.Range("C1").Value = "company"
$ws4.Range("D1").Value = "address"
$ws4.Range("E1").Value = "total"
$ws4.Range("F1").Value = "register_date"
$ws4.Range("G1").Value = "register_reason"
$ws4.Range("H1").Value = "property_category"
$ws4.Range("I1").Value = "category"
$ws4.Range("J1").Value = "date"
$ws4.Range("K1").Value = "legislator_name"
$ws4.Range("L1").Value = "legislator_id"
$ws4.Range("M1").Value = "source_file"
$ws4.Range("N1").Value = "index"

# Row 2
$ws4.Range("B2").Value = "顔清標"
$ws4.Range("C2").Value = "天台砂石股份有限公司"
$ws4.Range("D2").Value = "臺中市竹林里中山路紅竹巷58號1樓"
$ws4.Range("E2").Value = 1676000
$ws4.Range("F2").Value = "87年07月15R"
$ws4.Range("G2").Value = "合資公司"
$ws4.Range("H2").Value = "investment"
$ws4.Range("I2").Value = "normal"
$ws4.Range("J2").Value = "2012-04-10"
$ws4.Range("K2").Value = "顏清標"
$ws4.Range("L2").Value = 979
$ws4.Range("M2").Value = "tmp1b4d1"
$ws4.Range("N2").Value = 101

# Row 3
$ws4.Range("B3").Value = "顔清標"
$ws4.Range("C3").Value = "天馬瀝青股份有限公司"
$ws4.Range("D3").Value = "臺中市港埠路1段229號"
$ws4.Range("E3").Value = 3000000
$ws4.Range("F3").Value = "83年03月16日"
$ws4.Range("G3").Value = "合資公司"
$ws4.Range("H3").Value = "investment"
$ws4.Range("I3").Value = "normal"
$ws4.Range("J3").Value = "2012-04-10"
$ws4.Range("K3").Value = "顏清標"
$ws4.Range("L3").Value = 979
$ws4.Range("M3").Value = "tmp1b4d1"
$ws4.Range("N3").Value = 102

# Row 4
$ws4.Range("B4").Value = "顏清標"
$ws4.Range("C4").Value = "僑鴻建設股份有限公司"
$ws4.Range("D4").Value = "臺中市埔子里正義路1號"
$ws4.Range("E4").Value = 12000000
$ws4.Range("F4").Value = "85年07月13曰"
$ws4.Range("G4").Value = "合資公司"
$ws4.Range("H4").Value = "investment"
$ws4.Range("I4").Value = "normal"
$ws4.Range("J4").Value = "2012-04-10"
$ws4.Range("K4").Value = "顏清標"
$ws4.Range("L4").Value = 979
$ws4.Range("M4").Value = "tmp1b4d1"
$ws4.Range("N4").Value = 103

# Copy header style onto new header cells
$ws4.Range("G1").Copy()
$ws4.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Copy data row style onto new data cells (best-effort)
$ws4.Range("G2").Copy()
$ws4.Range("H2:N2").PasteSpecial(-4122)
$ws4.Range("G3").Copy()
$ws4.Range("H3:N3").PasteSpecial(-4122)
$ws4.Range("G4").Copy()
$ws4.Range("H4:N4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "edit complete"
